$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "71.133.48"
$ws.Range("E2").Value = "  +0.61%  "
Set-TextValue "D3" "3.863.33"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "700.25"
$ws.Range("E5").Value = "  +1.84%  "
Set-TextValue "D6" "173.76"
$ws.Range("E6").Value = "  +0.94%  "
Set-TextValue "D7" "3.863.06"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.525"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +0.63%  "
Set-TextValue "D11" "7.19"
$ws.Range("E11").Value = "  -3.40%  "
Set-TextValue "D12" "0.462"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  +4.66%  "
Set-TextValue "D14" "36.59"
$ws.Range("E14").Value = "  +1.13%  "
Set-TextValue "D15" "4.515.41"
$ws.Range("E15").Value = "  +1.23%  "
Set-TextValue "D16" "3.882.68"
$ws.Range("E16").Value = "  +1.65%  "
Set-TextValue "D17" "71.218.03"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "7.27"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "17.74"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +0.05%  "
Set-TextValue "D21" "11.21"
$ws.Range("E21").Value = "  -1.98%  "
Set-TextValue "D22" "496.77"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("E23").Value = "  +1.16%  "
Set-TextValue "D24" "85.37"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +0.50%  "
Set-TextValue "D26" "10.77"
$ws.Range("E26").Value = "  +3.42%  "
Set-TextValue "D27" "12.35"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D28" "4.008.99"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D29" "2.16"
$ws.Range("E29").Value = "  +0.42%  "
Set-TextValue "D30" "3.20"
$ws.Range("E30").Value = "  +6.82%  "
$ws.Range("E31").Value = "  -0.11%  "
Set-TextValue "D32" "7.68"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  -1.50%  "
Set-TextValue "D34" "29.71"
$ws.Range("E34").Value = "  -0.27%  "
Set-TextValue "D35" "0.182"
$ws.Range("E35").Value = "  -0.56%  "
Set-TextValue "D36" "9.30"
$ws.Range("E36").Value = "  +0.98%  "
Set-TextValue "D37" "3.815.91"
$ws.Range("E37").Value = "  +1.26%  "
Set-TextValue "D39" "0.105"
$ws.Range("E39").Value = "  +2.32%  "
Set-TextValue "D40" "2.42"
$ws.Range("E40").Value = "  +11.38%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D41" "3.41"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "6.07"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("E43").Value = "  +7.19%  "
$ws.Range("E44").Value = "  +0.07%  "
Set-TextValue "D46" "163.40"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "44.62"
$ws.Range("E48").Value = "  -3.45%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D49" "48.67"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  +1.14%  "
Set-TextValue "D51" "418.49"
$ws.Range("E51").Value = "  +4.66%  "
